$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the RMSE-based vote-prediction accuracy values for the
# "BASADO EN INSTANCIAS" table (color histogram HIST row + IncV3 row)
$ws.Range("B15").Value = 4.42
$ws.Range("C15").Value = 4.51
$ws.Range("D15").Value = 4.51
$ws.Range("E15").Value = 4.51

$ws.Range("B16").Value = 4.47
$ws.Range("C16").Value = 4.6100000000000003
$ws.Range("D16").Value = 4.6100000000000003
$ws.Range("E16").Value = 4.6100000000000003

$ws.Range("B19").Value = 4.08
$ws.Range("C19").Value = 4.24
$ws.Range("D19").Value = 4.34
$ws.Range("E19").Value = 4.4400000000000004

$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 4.2
$ws.Range("D20").Value = 4.47
$ws.Range("E20").Value = 4.5

# Update the saved view state (scroll position + active cell selection)
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("F19").Select() | Out-Null
